# DescData.xlsx: "unify the conception of DataNode, DataTable, Entity."
# The sheet that used to describe a generic "Property1" table is renamed to
# "DataNode", its two header rows are trimmed a touch, and the author leaves
# the selection sitting on D37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "Property1" to "DataNode"
$ws.Name = "DataNode"

# Header rows (row 1 and row 8) shrink from 28pt to 27pt
$ws.Rows(1).RowHeight = 27
$ws.Rows(8).RowHeight = 27

# Leave the active selection on D37 (sheet must be active for the
# selection to stick in the saved view state)
$ws.Activate()
[void]$ws.Range("D37").Select()
